# Apply unit change (Pa -> MPa style: divide by 1,000,000) and re-assign the
# Sector label ordering across every yearly worksheet in the AIC_Avg workbook.
#
# Previously: row 5 = "Offshore wind plants", row 6 = "Onshore wind plants",
#             row 7 = "Photovoltaic plants"
# Now:        row 5 = "Onshore wind plants",  row 6 = "Photovoltaic plants",
#             row 7 = "Offshore wind plants"
#
# The EU27+UK (column E) figures for those three rows are rotated along with
# the labels and rescaled (divided by 1,000,000) to reflect the new unit.

$wb = $excel.ActiveWorkbook

$newLabel5 = "Onshore wind plants"
$newLabel6 = "Photovoltaic plants"
$newLabel7 = "Offshore wind plants"

$count = $wb.Worksheets.Count
for ($i = 1; $i -le $count; $i++) {
    $ws = $wb.Worksheets.Item($i)

    $oldE5 = $ws.Range("E5").Value()
    $oldE6 = $ws.Range("E6").Value()
    $oldE7 = $ws.Range("E7").Value()

    $ws.Range("E5").Value = $oldE6 / 1000000
    $ws.Range("E6").Value = $oldE7 / 1000000
    $ws.Range("E7").Value = $oldE5 / 1000000

    $ws.Range("C5").Value = $newLabel5
    $ws.Range("C6").Value = $newLabel6
    $ws.Range("C7").Value = $newLabel7
}
